$wb = $excel.ActiveWorkbook

# Rename the third sheet from "macros" to "macro"
$wsMacro = $wb.Worksheets.Item("macros")
$wsMacro.Name = "macro"

# Update the selection on the "license" sheet (it currently holds the
# selection/active state) back to its default before moving focus away.
$wsLicense = $wb.Worksheets.Item("license")
$wsLicense.Range("B1").Select()

# Move the active selection on the "macro" sheet from B7 to B17 (it is
# split/frozen at row 1, so the selection lives in the bottom-left pane).
$wsMacro.Range("B17").Select()

# Make "macro" the active sheet/tab (was "license").
$wsMacro.Activate()
